$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header row: add N1 = 12, matching style of existing header cells (row 1) ---
$ws.Range("N1").Value = 12
$ws.Range("M1").Copy()
$ws.Range("N1").PasteSpecial(-4122)  # xlPasteFormats

# --- 2. Temporarily mark the whole data block (B2:N9) as Text so that
#        numeric-looking correlation values ("-0.06", "0.2", "0.0", ...)
#        are stored as text (shared strings) instead of being auto-converted
#        to numbers, matching the workbook's existing convention for this table. ---
$textHelper = $ws.Range("Z1")
$textHelper.NumberFormat = "@"
$textHelper.Copy()
$ws.Range("B2:N9").PasteSpecial(-4122)  # xlPasteFormats
$textHelper.Clear()

# --- 3. Write the new cross-correlation / Newey-West corrected values ---
$ws.Range("B2").Value = "-0.06"
$ws.Range("C2").Value = "-0.13"
$ws.Range("D2").Value = "-0.11"
$ws.Range("E2").Value = "-0.11"
$ws.Range("F2").Value = "-0.12"
$ws.Range("G2").Value = "-0.22*"
$ws.Range("H2").Value = "-0.18"
$ws.Range("I2").Value = "-0.25**"
$ws.Range("J2").Value = "-0.22*"
$ws.Range("K2").Value = "-0.22*"
$ws.Range("L2").Value = "-0.17"
$ws.Range("M2").Value = "-0.23*"
$ws.Range("N2").Value = "-0.07"
$ws.Range("B3").Value = "-0.02"
$ws.Range("C3").Value = "-0.09"
$ws.Range("D3").Value = "-0.08"
$ws.Range("E3").Value = "-0.07"
$ws.Range("F3").Value = "-0.04"
$ws.Range("G3").Value = "-0.16"
$ws.Range("H3").Value = "-0.15"
$ws.Range("I3").Value = "-0.22*"
$ws.Range("J3").Value = "-0.21*"
$ws.Range("K3").Value = "-0.2"
$ws.Range("L3").Value = "-0.13"
$ws.Range("M3").Value = "-0.21"
$ws.Range("N3").Value = "-0.05"
$ws.Range("B4").Value = "0.08"
$ws.Range("C4").Value = "0.14"
$ws.Range("D4").Value = "0.13"
$ws.Range("E4").Value = "0.07"
$ws.Range("F4").Value = "-0.02"
$ws.Range("G4").Value = "-0.14"
$ws.Range("H4").Value = "-0.18"
$ws.Range("I4").Value = "-0.31**"
$ws.Range("J4").Value = "-0.32**"
$ws.Range("K4").Value = "-0.28**"
$ws.Range("L4").Value = "-0.3**"
$ws.Range("M4").Value = "-0.26**"
$ws.Range("N4").Value = "-0.23*"
$ws.Range("B5").Value = "nan"
$ws.Range("C5").Value = "nan"
$ws.Range("D5").Value = "nan"
$ws.Range("E5").Value = "nan"
$ws.Range("F5").Value = "nan"
$ws.Range("G5").Value = "nan"
$ws.Range("H5").Value = "nan"
$ws.Range("I5").Value = "nan"
$ws.Range("J5").Value = "nan"
$ws.Range("K5").Value = "nan"
$ws.Range("L5").Value = "nan"
$ws.Range("M5").Value = "nan"
$ws.Range("N5").Value = "nan"
$ws.Range("B6").Value = "0.05"
$ws.Range("C6").Value = "0.11"
$ws.Range("D6").Value = "0.1"
$ws.Range("E6").Value = "0.03"
$ws.Range("F6").Value = "-0.1"
$ws.Range("G6").Value = "-0.06"
$ws.Range("H6").Value = "-0.12"
$ws.Range("I6").Value = "-0.2"
$ws.Range("J6").Value = "-0.24*"
$ws.Range("K6").Value = "-0.34***"
$ws.Range("L6").Value = "-0.39***"
$ws.Range("M6").Value = "-0.27**"
$ws.Range("N6").Value = "-0.33**"
$ws.Range("B7").Value = "0.11"
$ws.Range("C7").Value = "0.16"
$ws.Range("D7").Value = "0.15"
$ws.Range("E7").Value = "0.09"
$ws.Range("F7").Value = "-0.06"
$ws.Range("G7").Value = "-0.04"
$ws.Range("H7").Value = "-0.13"
$ws.Range("I7").Value = "-0.23*"
$ws.Range("J7").Value = "-0.29**"
$ws.Range("K7").Value = "-0.4***"
$ws.Range("L7").Value = "-0.44***"
$ws.Range("M7").Value = "-0.36***"
$ws.Range("N7").Value = "-0.38***"
$ws.Range("B8").Value = "0.12"
$ws.Range("C8").Value = "0.2"
$ws.Range("D8").Value = "0.31**"
$ws.Range("E8").Value = "0.25**"
$ws.Range("F8").Value = "0.06"
$ws.Range("G8").Value = "0.0"
$ws.Range("H8").Value = "-0.01"
$ws.Range("I8").Value = "-0.14"
$ws.Range("J8").Value = "-0.15"
$ws.Range("K8").Value = "-0.21"
$ws.Range("L8").Value = "-0.19"
$ws.Range("M8").Value = "-0.12"
$ws.Range("N8").Value = "-0.17"
$ws.Range("B9").Value = "-0.13"
$ws.Range("C9").Value = "-0.03"
$ws.Range("D9").Value = "-0.17"
$ws.Range("E9").Value = "-0.17"
$ws.Range("F9").Value = "-0.09"
$ws.Range("G9").Value = "-0.06"
$ws.Range("H9").Value = "-0.16"
$ws.Range("I9").Value = "-0.23*"
$ws.Range("J9").Value = "-0.2"
$ws.Range("K9").Value = "-0.18"
$ws.Range("L9").Value = "-0.26**"
$ws.Range("M9").Value = "-0.26**"
$ws.Range("N9").Value = "-0.27**"

# --- 4. Restore the plain (unstyled) look of the data cells, same as before the edit ---
$normalHelper = $ws.Range("Y1")
$normalHelper.Copy()
$ws.Range("B2:N9").PasteSpecial(-4122)  # xlPasteFormats
$normalHelper.Clear()
